$d = $word.ActiveDocument

# --- Simple text replacements -------------------------------------------------
# NB: "...Employee" is a textual prefix of the final "...Employee / Employer",
# so the plain-"Employee" paragraph must be rewritten *before* the plain-
# "Employer" one, otherwise the later global Find would also re-match (and
# re-append to) the text this step just wrote.

$d.Content.Find.Execute(
    "Predicate (Role): AggregationSubjectKind. Employee", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Predicate (Role): AggregationSubjectKind. Employee / Employer", 2)

$d.Content.Find.Execute(
    "Predicate (Role): AggregationSubjectKind. Employer", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Predicate (Role): AggregationSubjectKind. Employee / Employer", 2)

$d.Content.Find.Execute(
    "Object (Occurrence): AggregationSubject. anEmployer", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Object (Occurrence): AggregationSubject. anEmployee / anEmployer", 2)

$d.Content.Find.Execute(
    "Subject (Context): AggregationSubject. anEmployee", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Subject (Context): AggregationSubject. anEmployee / anEmployer", 2)

# --- TODO. paragraph -> "C: Transform Context. Mapping State." plus new
#     trailing paragraphs describing the flatMap composition -----------------

$d.Content.Find.Execute(
    "TODO.", $true, $false, $false, $false, $false, $true, 1, $false,
    "C: Transform Context. Mapping State.", 2)

$rng = $d.Content
$rng.Find.Execute("C: Transform Context. Mapping State.")
$curPara = $rng.Paragraphs(1)

$newTexts = @(
    "P: Transforms, Mappings.",
    "S, O: Monads, Functors.",
    "",
    "S.flatMap(P) : O;",
    "O.flatMap(P) : S;"
)

foreach ($t in $newTexts) {
    $curPara.Range.InsertParagraphAfter()
    $curPara = $curPara.Next()
    if ($t -ne "") {
        $curPara.Range.Text = $t
    }
}
